# Generate Report for Handoff
# Update the "Ready for handoff" rows: set Priority to "ht" and bump the
# handoff timestamps on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 10, 11, 12, 13)

foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-26 06:20:52"
    $wsZhCn.Range("H$r").Value     = "2016-08-26 06:20:47"
    $wsDeDe.Range("H$r").Value     = "2016-08-26 06:20:52"

    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"
}
